# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the Price (D) cells that receive new values,
# so numeric-looking strings (e.g. '1.00', '241.48') are preserved verbatim
# instead of being coerced into floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply new cell values
$ws.Range("D2").Value = "36.471.96"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "2.061.33"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "241.48"
$ws.Range("E5").Value = "  -2.92%  "
$ws.Range("D6").Value = "0.659"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "52.32"
$ws.Range("E8").Value = "  -8.48%  "
$ws.Range("D9").Value = "58.57"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").Value = "0.358"
$ws.Range("E10").Value = "  -7.54%  "
$ws.Range("D11").Value = "0.0745"
$ws.Range("E11").Value = "  -5.13%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "0.882"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").Value = "14.58"
$ws.Range("E14").Value = "  -10.40%  "
$ws.Range("D15").Value = "2.375.73"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "5.39"
$ws.Range("E16").Value = "  -6.91%  "
$ws.Range("D17").Value = "2.097.48"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("D18").Value = "36.514.95"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").Value = "16.44"
$ws.Range("E19").Value = "  -12.10%  "
$ws.Range("D20").Value = "71.37"
$ws.Range("E20").Value = "  -4.73%  "
$ws.Range("D21").Value = "0.0₃0855"
$ws.Range("E21").Value = "  -5.14%  "
$ws.Range("D22").Value = "5.25"
$ws.Range("E22").Value = "  -4.63%  "
$ws.Range("D23").Value = "236.29"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "2.37"
$ws.Range("E25").Value = "  -4.93%  "
$ws.Range("D26").Value = "9.42"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("D27").Value = "2.12"
$ws.Range("E27").Value = "  -3.50%  "
$ws.Range("D28").Value = "164.03"
$ws.Range("E28").Value = "  -3.50%  "
$ws.Range("D29").Value = "20.28"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "0.121"
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").Value = "5.04"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("D32").Value = "1.14"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("D33").Value = "4.56"
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").Value = "0.0589"
$ws.Range("E34").Value = "  -5.84%  "
$ws.Range("D35").Value = "2.35"
$ws.Range("E35").Value = "  +2.97%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").Value = "1.83"
$ws.Range("E37").Value = "  +3.37%  "
$ws.Range("E38").Value = "  -8.95%  "
$ws.Range("D39").Value = "1.24"
$ws.Range("E39").Value = "  -8.09%  "
$ws.Range("D40").Value = "4.83"
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0216"
$ws.Range("E41").Value = "  -3.66%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "1.13"
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "2.91"
$ws.Range("E43").Value = "  -6.11%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.0938"
$ws.Range("E44").Value = "  -7.81%  "
$ws.Range("D45").Value = "93.58"
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "7.52"
$ws.Range("E46").Value = "  +9.72%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.376.13"
$ws.Range("E47").Value = "  +7.78%  "
$ws.Range("D48").Value = "15.32"
$ws.Range("E48").Value = "  -13.59%  "
$ws.Range("D49").Value = "2.34"
$ws.Range("E49").Value = "  -4.45%  "
$ws.Range("D50").Value = "2.87"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "2.264.02"
$ws.Range("E51").Value = "  +0.78%  "
